# Weekly refresh of the Cilantro price sheet:
#  - insert 2 new rows (latest week, Fecha 44783) right before the current
#    row 641, shifting every subsequent row down by 2
#  - append 2 rows at the end of the table that repeat the last data pair
#    (Fecha 44572) which fell out of the shifted range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the date column's number format so newly-created cells keep the
# same display (YYYY-MM-DD HH:MM:SS) as the rest of column D.
$dateFormat = $ws.Range("D2").NumberFormat

function Set-CilantroRow {
    param(
        $Row,
        $Fecha,
        $Volumen,
        $PrecioMin,
        $PrecioMax,
        $PrecioProm,
        $Unidad,
        $PrecioKg,
        $KgUnidades
    )

    $ws.Range("A$Row").Value = 9
    $ws.Range("B$Row").Value = "Vega Central Mapocho de Santiago"
    $ws.Range("C$Row").Value = "Metropolitana"

    $ws.Range("D$Row").NumberFormat = $dateFormat
    $ws.Range("D$Row").Value = $Fecha

    $ws.Range("E$Row").Value = 13
    $ws.Range("F$Row").Value = 100112040
    $ws.Range("G$Row").Value = "Cilantro"
    $ws.Range("H$Row").Value = "Sin especificar"
    $ws.Range("I$Row").Value = "Primera"
    $ws.Range("J$Row").Value = $Volumen
    $ws.Range("K$Row").Value = $PrecioMin
    $ws.Range("L$Row").Value = $PrecioMax
    $ws.Range("M$Row").Value = $PrecioProm
    $ws.Range("N$Row").Value = $Unidad
    $ws.Range("O$Row").Value = "Región Metropolitana"
    $ws.Range("P$Row").Value = $PrecioKg
    $ws.Range("Q$Row").Value = $KgUnidades
    $ws.Range("R$Row").Value = "Hortaliza"
}

# Insert two fresh rows before row 641 (everything from 641 on shifts down by 2).
$ws.Range("A641:A642").EntireRow.Insert()

# New week's data (Fecha = 44783), landing in the freshly inserted rows.
Set-CilantroRow 641 44783 52  11000 11000 11000 "`$/caja 36 atados"   306  36
Set-CilantroRow 642 44783 160 15000 16000 15500 "`$/docena de atados" 5167 3

# The two rows that used to close the table (old 692/693, Fecha = 44572) are
# re-appended at the new end of the table (694/695).
Set-CilantroRow 694 44572 43  8000  8000  8000  "`$/caja 36 atados"   222  36
Set-CilantroRow 695 44572 106 16000 18000 17000 "`$/docena de atados" 5667 3
